$d = $word.ActiveDocument
$p2 = $d.Paragraphs(2)
$rng = $p2.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00BB345C" w:rsidRPr="008547DF" w:rsidRDefault="00F05369" w:rsidP="00F05369"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r w:rsidRPr="00F05369"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Será usado neste projeto o ciclo de vida de produção baseado em </w:t></w:r><w:r w:rsidR="008547DF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>metodologias ágeis</w:t></w:r><w:r w:rsidRPr="00F05369"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, o projeto necessita de flexibilidade para mudanças futuras como novas funcionalidades, remoção de outras e melhorias, além da necessidade da entrega continua de partes do projeto para o cliente, desta forma a melhor metodologia a se adotar é a</w:t></w:r><w:r w:rsidR="008547DF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> de metodologias ágeis</w:t></w:r><w:r w:rsidRPr="00F05369"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, com ela é possível definir o escopo do projeto inicialmente mas não impede de adicionar ou remover </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>funcionalidades</w:t></w:r><w:r w:rsidRPr="00F05369"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> além de ser possível dividir o projeto em pequenas partes, podendo ser entregues como o cliente exigiu.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">O método SCRUM pode ser aplicado ao projeto pois ele faz parte das metodologias ágeis, como foi citado anteriormente o cliente necessita de uma flexibilidade no projeto e entrega continua de resultados, com o SCRUM isso é possível de maneira eficiente e fácil, é </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>possível</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> fazer alterações no projeto em </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>qualquer momento sem nenhum problema.</w:t></w:r></w:p>'

$rng.InsertXML($xml)

Write-Output "Paragraph replaced; new paragraph count:"
Write-Output $d.Paragraphs.Count
